# Add two new columns, I ("I0") and J ("IF"), to the existing results
# table on the active sheet, mirroring the style of the existing header
# row (H is the last existing column, bold / bordered / centered) and
# filling in the per-row numeric values for rows 2-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the header formatting (bold font, thin border box, centered
# alignment) from the existing "IP" header cell (H1) onto the two new
# header cells so they match the look of the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row values for the new I0 / IF columns.
$ijValues = @{
    2  = @(8, 8)
    3  = @(9, 9)
    4  = @(9, 9)
    5  = @(9, 9)
    6  = @(8, 8)
    7  = @(8, 8)
    8  = @(8, 9)
    9  = @(9, 9)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(8, 9)
    13 = @(8, 9)
    14 = @(8, 9)
    15 = @(8, 8)
    16 = @(6, 6)
    17 = @(9, 9)
    18 = @(8, 8)
    19 = @(8, 8)
    20 = @(6, 6)
    21 = @(6, 6)
    22 = @(7, 7)
    23 = @(6, 6)
    24 = @(3, 3)
}

foreach ($row in $ijValues.Keys) {
    $vals = $ijValues[$row]
    $ws.Cells.Item($row, 9).Value  = $vals[0]   # column I
    $ws.Cells.Item($row, 10).Value = $vals[1]   # column J
}
